$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HeroSkill")

# Balance change: "瞄准射击" (Aimed Shot) skill description updated from
# dealing magic damage to a single target (double vs heroes) to instead
# dealing magic damage to the enemy king tower.
$ws.Range("C7").Value = "对敌王塔造成魔法伤害"

# Reflect the cell selection left behind by the edit.
$ws.Range("C7").Select()
